# case_follow_up.xlsx - "fix several issues with updating CSV files when a
# form has a repository": add a new "Project" lookup question
# (select_one_from_file generated_case.csv / me_project) to the survey
# sheet, right after the "distrito" question and before the existing
# "households.csv" lookup (hid).

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

# Insert a new blank row at row 6 (shifts rows 6-10 down to 7-11), which
# also pushes the sheet's used range from A1:F10 to A1:F11. The inserted
# row inherits the formatting of the row above it, matching styles
# s="10"/s="10"/s="11"/s="11"/s="12" for columns A-E.
$survey.Rows("6:6").Insert()

# Populate the new XLSForm question: a select_one_from_file on
# generated_case.csv, named me_project, labeled "Project". The
# required / choice_filter / calculation columns (D-F) stay blank for
# this row.
$survey.Range("A6").Value = "select_one_from_file generated_case.csv"
$survey.Range("B6").Value = "me_project"
$survey.Range("C6").Value = "Project"

# Mirror the editor's cell selection left behind on the other two sheets
# after making this change (cosmetic view-state only).
$choices = $wb.Worksheets.Item("choices")
$choices.Range("C4").Select()

$settings = $wb.Worksheets.Item("settings")
$settings.Range("B3").Select()

# Leave the workbook focused back on the survey sheet, with the newly
# inserted row selected, same as after the original edit.
$survey.Activate()
$survey.Range("A6:C6").Select()
